# Base API local Grid Cloud
#
# - TestCases!B3 flips from "Y" to "N" (adds a new shared string "N").
# - The active/selected sheet moves from "TestCases" to "TestData":
#     * TestData becomes the tab-selected sheet (and workbook's activeTab).
#     * TestData's scroll/selection moves off B8 onto G6.
#     * TestCases loses its tabSelected flag.

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("TestCases")
$wsData  = $wb.Worksheets.Item("TestData")

# Flip the Runmode flag for the OpenAccountTest row from Y to N.
$wsCases.Range("B3").Value = "N"

# Make TestData the active sheet and move its selection to G6.
$wsData.Activate()
$wsData.Range("G6").Select()
